$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 1
$ws.Range("A26").Value = 2
$ws.Range("A27").Value = 2
$ws.Range("A28").Value = 2
$ws.Range("A33").Value = 4

$ws.Range("A29").Select()
